$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.657.56"
$ws.Range("E2").Value = "  +0.00%  "

$ws.Range("D3").Value = "'1.597.25"
$ws.Range("E3").Value = "  -0.01%  "

$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").Value = "'211.39"
$ws.Range("E5").Value = "  +0.03%  "

$ws.Range("E6").Value = "  +0.68%  "

$ws.Range("E7").Value = "  +0.19%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").Value = "  +0.91%  "

$ws.Range("D10").Value = "'19.56"
$ws.Range("E10").Value = "  -0.64%  "

$ws.Range("D11").Value = "'0.0840"
$ws.Range("E11").Value = "  +0.33%  "

$ws.Range("D12").Value = "'1.821.76"
$ws.Range("E12").Value = "  +0.05%  "

$ws.Range("D13").Value = "'1.585.94"
$ws.Range("E13").Value = "  -1.24%  "

$ws.Range("E14").Value = "  -0.18%  "

$ws.Range("D15").Value = "'0.523"
$ws.Range("E15").Value = "  +0.34%  "

$ws.Range("D16").Value = "'65.07"
$ws.Range("E16").Value = "  +0.30%  "

$ws.Range("D17").Value = "'26.650.69"
$ws.Range("E17").Value = "  -0.01%  "

$ws.Range("D18").Value = "'0.0₃0738"
$ws.Range("E18").Value = "  +1.40%  "

$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").Value = "'1.00"
$ws.Range("E19").Value = "  +0.17%  "

$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'209.06"
$ws.Range("E20").Value = "  -0.29%  "

$ws.Range("D21").Value = "'7.01"
$ws.Range("E21").Value = "  +3.46%  "

$ws.Range("E22").Value = "  +0.36%  "

$ws.Range("E23").Value = "  +1.46%  "

$ws.Range("D24").Value = "'8.98"
$ws.Range("E24").Value = "  +0.71%  "

$ws.Range("D25").Value = "'144.18"
$ws.Range("E25").Value = "  -1.41%  "

$ws.Range("E26").Value = "  +0.24%  "

$ws.Range("E27").Value = "  -0.93%  "

$ws.Range("E28").Value = "  -0.70%  "

$ws.Range("D29").Value = "'15.28"
$ws.Range("E29").Value = "  -0.07%  "

$ws.Range("D30").Value = "'0.0516"
$ws.Range("E30").Value = "  +2.37%  "

$ws.Range("D31").Value = "'1.16"
$ws.Range("E31").Value = "  +0.50%  "

$ws.Range("E32").Value = "  +0.48%  "

$ws.Range("E33").Value = "  +1.42%  "

$ws.Range("D34").Value = "'1.286.92"
$ws.Range("E34").Value = "  -1.10%  "

$ws.Range("D35").Value = "'0.617"
$ws.Range("E35").Value = "  -7.33%  "

$ws.Range("E36").Value = "  +0.48%  "

$ws.Range("E37").Value = "  +0.43%  "

$ws.Range("E38").Value = "  -0.59%  "

$ws.Range("E39").Value = "  -0.87%  "

$ws.Range("E40").Value = "  +18.04%  "

$ws.Range("D41").Value = "'5.49"
$ws.Range("E41").Value = "  +2.13%  "

$ws.Range("E42").Value = "  +0.10%  "

$ws.Range("D43").Value = "'0.782"
$ws.Range("E43").Value = "  -0.78%  "

$ws.Range("D44").Value = "'63.57"
$ws.Range("E44").Value = "  -0.39%  "

$ws.Range("D45").Value = "'1.734.35"
$ws.Range("E45").Value = "  +0.02%  "

$ws.Range("D46").Value = "'90.65"
$ws.Range("E46").Value = "  +0.46%  "

$ws.Range("D47").Value = "'1.56"
$ws.Range("E47").Value = "  -3.52%  "

$ws.Range("D48").Value = "'0.102"
$ws.Range("E48").Value = "  +1.49%  "

$ws.Range("E49").Value = "  +0.92%  "

$ws.Range("E50").Value = "  +0.13%  "

$ws.Range("D51").Value = "'7.40"
$ws.Range("E51").Value = "  -1.03%  "
